$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "E"=3; "F"=1; "G"=96.17831799999999; "H"=288.534954; "I"=0.7237598617297997; "J"=0.7237598617297996; "K"=3; "L"=1; "M"=117.044563; "N"=351.133689; "O"=0.3245365645427815; "P"=0.3245365645427815; "Q"=11257.14920038503; "R"=101314.3428034653; "S"=0.2348865390797477; "T"=0.2348865390797477 }
    3 = @{ "E"=3; "F"=1; "G"=96.17831799999999; "H"=288.534954; "I"=0.7237598617297997; "J"=0.7237598617297996; "K"=3; "L"=1; "M"=101.5800373333333; "N"=304.740112; "O"=0.281657135515876; "P"=0.281657135515876; "Q"=9769.797133097205; "R"=87928.17419787485; "S"=0.2038521294561819; "T"=0.2038521294561818 }
    4 = @{ "E"=3; "F"=1; "G"=96.17831799999999; "H"=288.534954; "I"=0.7237598617297997; "J"=0.7237598617297996; "K"=3; "L"=1; "M"=142.0267893333333; "N"=426.080368; "O"=0.3938062999413425; "P"=0.3938062999413425; "Q"=13659.89770902034; "R"=122939.0793811831; "S"=0.2850211931938701; "T"=0.28502119319387 }
    5 = @{ "E"=3; "F"=1; "G"=13.23504133333333; "H"=39.705124; "I"=0.09959616558694152; "J"=0.0995961655869415; "K"=3; "L"=1; "M"=117.044563; "N"=351.133689; "O"=0.3245365645427815; "P"=0.3245365645427815; "Q"=1549.089629146937; "R"=13941.80666232243; "S"=0.03232259742122; "T"=0.03232259742121999 }
    6 = @{ "E"=3; "F"=1; "G"=13.23504133333333; "H"=39.705124; "I"=0.09959616558694152; "J"=0.0995961655869415; "K"=3; "L"=1; "M"=101.5800373333333; "N"=304.740112; "O"=0.281657135515876; "P"=0.281657135515876; "Q"=1344.41599274821; "R"=12099.74393473389; "S"=0.02805197070758281; "T"=0.0280519707075828 }
    7 = @{ "E"=3; "F"=1; "G"=13.23504133333333; "H"=39.705124; "I"=0.09959616558694152; "J"=0.0995961655869415; "K"=3; "L"=1; "M"=142.0267893333333; "N"=426.080368; "O"=0.3938062999413425; "P"=0.3938062999413425; "Q"=1879.730427267292; "R"=16917.57384540563; "S"=0.03922159745813871; "T"=0.0392215974581387 }
    8 = @{ "E"=3; "F"=1; "G"=23.47369766666667; "H"=70.421093; "I"=0.1766439726832589; "J"=0.1766439726832589; "K"=3; "L"=1; "M"=117.044563; "N"=351.133689; "O"=0.3245365645427815; "P"=0.3245365645427815; "Q"=2747.468685389119; "R"=24727.21816850208; "S"=0.05732742804181379; "T"=0.05732742804181378 }
    9 = @{ "E"=3; "F"=1; "G"=23.47369766666667; "H"=70.421093; "I"=0.1766439726832589; "J"=0.1766439726832589; "K"=3; "L"=1; "M"=101.5800373333333; "N"=304.740112; "O"=0.281657135515876; "P"=0.281657135515876; "Q"=2384.45908533138; "R"=21460.13176798242; "S"=0.04975303535211136; "T"=0.04975303535211134 }
    10 = @{ "E"=3; "F"=1; "G"=23.47369766666667; "H"=70.421093; "I"=0.1766439726832589; "J"=0.1766439726832589; "K"=3; "L"=1; "M"=142.0267893333333; "N"=426.080368; "O"=0.3938062999413425; "P"=0.3938062999413425; "Q"=3333.893913378025; "R"=30005.04522040223; "S"=0.06956350928933379; "T"=0.06956350928933377 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($c in $rowVals.Keys) {
        $ws.Range("$c$r").Value = $rowVals[$c]
    }
}

Write-Output "Done updating cells"